# Update "想去人数" (F column) values across the 4 sheets to match the
# regenerated data snapshot (gh-pages output at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3321
$ws.Range("F6").Value = 7747
$ws.Range("F9").Value = 733
$ws.Range("F14").Value = 172
$ws.Range("F15").Value = 1760
$ws.Range("F16").Value = 371
$ws.Range("F17").Value = 91
$ws.Range("F18").Value = 2371
$ws.Range("F20").Value = 1034
$ws.Range("F21").Value = 1031
$ws.Range("F22").Value = 1046
$ws.Range("F23").Value = 6348
$ws.Range("F24").Value = 6495
$ws.Range("F25").Value = 396
$ws.Range("F27").Value = 1087
$ws.Range("F30").Value = 521
$ws.Range("F31").Value = 1075
$ws.Range("F33").Value = 241
$ws.Range("F34").Value = 241
$ws.Range("F37").Value = 84
$ws.Range("F38").Value = 601
$ws.Range("F39").Value = 416
$ws.Range("F41").Value = 1250
$ws.Range("F42").Value = 3252
$ws.Range("F44").Value = 716
$ws.Range("F46").Value = 44
$ws.Range("F49").Value = 470
$ws.Range("F50").Value = 62

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 82
$ws.Range("F25").Value = 3
$ws.Range("F26").Value = 6609
$ws.Range("F42").Value = 23

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1996
$ws.Range("F5").Value = 1321
$ws.Range("F7").Value = 555
$ws.Range("F8").Value = 2148
$ws.Range("F9").Value = 8932
$ws.Range("F10").Value = 1050

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3321
$ws.Range("F4").Value = 1996
$ws.Range("F5").Value = 1321
$ws.Range("F6").Value = 555
$ws.Range("F7").Value = 2148
$ws.Range("F9").Value = 1050
$ws.Range("F14").Value = 172
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 2371
$ws.Range("F19").Value = 1034
$ws.Range("F20").Value = 1031
$ws.Range("F21").Value = 1046
$ws.Range("F22").Value = 6348
$ws.Range("F23").Value = 6496
$ws.Range("F24").Value = 396
$ws.Range("F26").Value = 1087
$ws.Range("F29").Value = 521
$ws.Range("F30").Value = 1075
$ws.Range("F31").Value = 241
$ws.Range("F32").Value = 241
$ws.Range("F35").Value = 84
$ws.Range("F36").Value = 601
$ws.Range("F37").Value = 416
$ws.Range("F40").Value = 3252
$ws.Range("F41").Value = 716
$ws.Range("F44").Value = 6609
